# Update "想去人数" (column F) values for the 展览 sheet and the
# combined 全部类型 sheet to reflect newly refreshed scrape data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row -> new value on the "展览" sheet
$sheet1Updates = @{
    4  = 1755
    7  = 1141
    12 = 3123
    13 = 667
    14 = 1809
    15 = 1818
    16 = 878
    19 = 1502
    20 = 299
    23 = 1271
    26 = 166
    27 = 5828
    28 = 5321
    31 = 1692
    33 = 202
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Row -> new value on the "全部类型" sheet
$sheet4Updates = @{
    7  = 1755
    10 = 1141
    16 = 3123
    17 = 667
    18 = 1809
    19 = 1818
    20 = 878
    23 = 1502
    24 = 299
    29 = 1271
    32 = 166
    33 = 5828
    34 = 5321
    37 = 1692
    41 = 202
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
